$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: new columns I and J, matching the existing header style (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-5 for columns I and J (plain numeric cells, no special style)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9
